$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 2
    11 = 1
    12 = 1
    13 = 2
    14 = 0
    15 = 2
    16 = 3
    17 = 3
    18 = 1
    19 = 0
    20 = 3
    21 = 2
    22 = 1
    23 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
